$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format on price cells that would otherwise be auto-converted to numbers
$textCells = @("D4", "D5", "D6", "D9", "D10", "D12", "D13", "D14", "D18", "D20", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D35", "D38", "D39", "D42", "D45", "D46", "D47", "D48", "D49", "D21", "D22", "D40", "D41", "D50", "D51")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Update Price (D) and Volume(1h) (E) columns for rows with same coin (no reordering)
$ws.Range("D2").Value = "71.094.72"
$ws.Range("E2").Value = "  +0.07%  "
$ws.Range("D3").Value = "3.845.02"
$ws.Range("E3").Value = "  +0.65%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "704.00"
$ws.Range("E5").Value = "  +0.37%  "
$ws.Range("D6").Value = "172.90"
$ws.Range("E6").Value = "  -0.16%  "
$ws.Range("D7").Value = "3.841.52"
$ws.Range("E7").Value = "  +0.60%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").Value = "0.526"
$ws.Range("E9").Value = "  -0.55%  "
$ws.Range("D10").Value = "0.163"
$ws.Range("E10").Value = "  +0.51%  "
$ws.Range("E11").Value = "  +0.71%  "
$ws.Range("D12").Value = "0.459"
$ws.Range("E12").Value = "  -0.14%  "
$ws.Range("D13").Value = "0.0000258"
$ws.Range("E13").Value = "  +0.33%  "
$ws.Range("D14").Value = "36.63"
$ws.Range("E14").Value = "  +0.69%  "
$ws.Range("D15").Value = "4.495.94"
$ws.Range("E15").Value = "  +0.81%  "
$ws.Range("D16").Value = "3.862.85"
$ws.Range("E16").Value = "  +1.38%  "
$ws.Range("D17").Value = "71.132.24"
$ws.Range("E17").Value = "  +0.25%  "
$ws.Range("D18").Value = "7.20"
$ws.Range("E18").Value = "  -0.06%  "
$ws.Range("E19").Value = "  +0.67%  "
$ws.Range("D20").Value = "17.32"
$ws.Range("E20").Value = "  -2.93%  "
$ws.Range("D23").Value = "0.723"
$ws.Range("E23").Value = "  +1.22%  "
$ws.Range("D24").Value = "85.08"
$ws.Range("E24").Value = "  +1.41%  "
$ws.Range("D25").Value = "0.0000147"
$ws.Range("E25").Value = "  +2.15%  "
$ws.Range("D26").Value = "10.62"
$ws.Range("E26").Value = "  +1.38%  "
$ws.Range("D27").Value = "12.19"
$ws.Range("E27").Value = "  -1.62%  "
$ws.Range("D28").Value = "2.10"
$ws.Range("E28").Value = "  -3.13%  "
$ws.Range("D29").Value = "3.18"
$ws.Range("E29").Value = "  +2.18%  "
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").Value = "  -0.05%  "
$ws.Range("D31").Value = "7.49"
$ws.Range("E31").Value = "  -0.78%  "
$ws.Range("D32").Value = "2.26"
$ws.Range("E32").Value = "  -1.97%  "
$ws.Range("D33").Value = "29.41"
$ws.Range("E33").Value = "  -0.60%  "
$ws.Range("E34").Value = "  -2.95%  "
$ws.Range("D35").Value = "9.19"
$ws.Range("E35").Value = "  -0.70%  "
$ws.Range("D36").Value = "3.805.22"
$ws.Range("E36").Value = "  +1.08%  "
$ws.Range("E37").Value = "  -0.07%  "
$ws.Range("D38").Value = "0.103"
$ws.Range("E38").Value = "  -0.27%  "
$ws.Range("D39").Value = "2.36"
$ws.Range("E39").Value = "  +5.94%  "
$ws.Range("D42").Value = "3.36"
$ws.Range("E42").Value = "  -1.76%  "
$ws.Range("E43").Value = "  +0.01%  "
$ws.Range("E44").Value = "  +0.21%  "
$ws.Range("D45").Value = "0.000317"
$ws.Range("E45").Value = "  -2.12%  "
$ws.Range("D46").Value = "163.30"
$ws.Range("E46").Value = "  +0.18%  "
$ws.Range("D47").Value = "48.65"
$ws.Range("E47").Value = "  -0.45%  "
$ws.Range("D48").Value = "415.63"
$ws.Range("E48").Value = "  +1.73%  "
$ws.Range("D49").Value = "1.38"
$ws.Range("E49").Value = "  -0.43%  "

# Update rows that swapped coin ranking order (Coin, Link, Price, Volume)
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").Value = "10.71"
$ws.Range("E21").Value = "  -4.57%  "
$ws.Range("B22").Value = "BitcoinCash"
$ws.Range("C22").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D22").Value = "493.63"
$ws.Range("E22").Value = "  +2.67%  "
$ws.Range("B40").Value = "Mantle"
$ws.Range("C40").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D40").Value = "1.04"
$ws.Range("E40").Value = "  +5.48%  "
$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").Value = "6.02"
$ws.Range("E41").Value = "  -0.31%  "
$ws.Range("B50").Value = "Cosmos"
$ws.Range("C50").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D50").Value = "8.60"
$ws.Range("E50").Value = "  +0.28%  "
$ws.Range("B51").Value = "TheGraph"
$ws.Range("C51").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D51").Value = "0.297"
$ws.Range("E51").Value = "  -1.33%  "
